$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

# --- Update rows 150-159: mark predictions as completed with results ---
$updates = @(
    @{ Row = 150; L = "Completed"; M = "Draw";     N = "Fallo";   O = -2.9;  P = -100; },
    @{ Row = 151; L = "Completed"; M = "Draw";     N = "Fallo";   O = -2.2;  P = -100; },
    @{ Row = 152; L = "Completed"; M = "Away Win"; N = "Acierto"; O = 1.36;  P = 91;   },
    @{ Row = 153; L = "Completed"; M = "Draw";     N = "Fallo";   O = -1.9;  P = -100; },
    @{ Row = 154; L = "Completed"; M = "Away Win"; N = "Fallo";   O = -2;    P = -100; },
    @{ Row = 155; L = "Completed"; M = "Away Win"; N = "Acierto"; O = 1.6;   P = 80;   },
    @{ Row = 156; L = "Completed"; M = "Home Win"; N = "Acierto"; O = 1.3;   P = 45;   },
    @{ Row = 157; L = "Completed"; M = "Home Win"; N = "Acierto"; O = 1.45;  P = 50;   },
    @{ Row = 158; L = "Completed"; M = "Home Win"; N = "Acierto"; O = 1.69;  P = 65;   },
    @{ Row = 159; L = "Completed"; M = "Home Win"; N = "Acierto"; O = 0.88;  P = 110;  }
)

$timestamp = "2025-09-18 04:27:14"

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 12).Value = $u.L
    $ws.Cells.Item($r, 13).Value = $u.M
    $ws.Cells.Item($r, 14).Value = $u.N
    $ws.Cells.Item($r, 15).Value = $u.O
    $ws.Cells.Item($r, 16).Value = $u.P
    $ws.Cells.Item($r, 17).Value = $timestamp
}

# --- Append new rows 164-169 ---
$newRows = @(
    @{ Row = 164; A = "2025-09-19"; B = "Eredivisie";           C = "Sparta Rotterdam"; D = "Twente";            E = "Away Win"; F = "45.05%"; G = 2.35; H = "4.81%";  I = 0.3; J = 0.004344782829934631; K = 0.04344782829934631 },
    @{ Row = 165; A = "2025-09-19"; B = "Bundesliga";           C = "VfB Stuttgart";    D = "FC St. Pauli";      E = "Home Win"; F = "75.61%"; G = 1.73; H = "29.49%"; I = 2.5; J = 0.04218802910321603;  K = 0.4218802910321602  },
    @{ Row = 166; A = "2025-09-19"; B = "Jupiler Pro League";   C = "Gent";             D = "Dender";            E = "Home Win"; F = "71.79%"; G = 1.73; H = "22.96%"; I = 1.9; J = 0.03314905704217487;  K = 0.3314905704217487  },
    @{ Row = 167; A = "2025-09-19"; B = "La Liga";              C = "Real Betis";       D = "Real Sociedad";     E = "Home Win"; F = "62.37%"; G = 2;    H = "23.49%"; I = 1.4; J = 0.02473234525645181;  K = 0.247323452564518   },
    @{ Row = 168; A = "2025-09-19"; B = "Primeira Liga";        C = "Rio Ave";          D = "FC Porto";          E = "Away Win"; F = "90.96%"; G = 1.42; H = "27.87%"; I = 2.9; J = 0.05;                  K = 0.6943610353896599  },
    @{ Row = 169; A = "2025-09-19"; B = "Liga de Expansión MX"; C = "Tapatío";          D = "Correcaminos Uat";  E = "Home Win"; F = "71.94%"; G = 1.65; H = "17.51%"; I = 1.7; J = 0.02876458510447492;  K = 0.2876458510447492  }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    # Columns A, F and H hold text that looks like a date / percentage
    # ("2025-09-19", "45.05%", ...) in the source data - force text format
    # so Excel doesn't auto-convert them into a date serial / numeric value.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $nr.A
    $ws.Cells.Item($r, 2).Value = $nr.B
    $ws.Cells.Item($r, 3).Value = $nr.C
    $ws.Cells.Item($r, 4).Value = $nr.D
    $ws.Cells.Item($r, 5).Value = $nr.E
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $nr.F
    $ws.Cells.Item($r, 7).Value = $nr.G
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = $nr.H
    $ws.Cells.Item($r, 9).Value = $nr.I
    $ws.Cells.Item($r, 10).Value = $nr.J
    $ws.Cells.Item($r, 11).Value = $nr.K
    $ws.Cells.Item($r, 12).Value = "Pending"

    # Result columns (M-Q) are still blank for these newly-added, not-yet-played
    # fixtures - match the "Pending" rows already present earlier in the sheet.
    $ws.Cells.Item($r, 13).Value = ""
    $ws.Cells.Item($r, 14).Value = ""
    $ws.Cells.Item($r, 15).Value = ""
    $ws.Cells.Item($r, 16).Value = ""
    $ws.Cells.Item($r, 17).Value = ""
}
